$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets list the same events in rows 2-5 and
# need their "想去人数" (F column) counts refreshed to the newly scraped
# values.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 142
    $ws.Range("F3").Value = 217
    $ws.Range("F4").Value = 3695
    $ws.Range("F5").Value = 382
}
